$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "optimization_parameters": restructure to the new "beta" format.
#   - drop the duplicated "value" header cells in C1:F1
#   - rename "Model" -> "production_function"
#   - insert a new "L_curve" parameter row right after "production_function"
#   - drop the obsolete "Deletion" row
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Remove the extra duplicate "value" cells in row 1 (C1:F1) -> spans shrink to 1:5
$ws.Range("C1:F1").ClearContents()

# Insert a new blank row at 9 (current rows 9-17 shift down to 10-18)
$ws.Rows.Item(9).Insert()

# Row 8: "Model" becomes "production_function" (value stays "Sigmoid")
$ws.Range("A8").Value = "production_function"

# New row 9: "L_curve" parameter, defaulting to 0, scientific-notation formatted
# like the other optimizer knobs above it.
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (originally row 16) is now at row 17 after the insert
# above; drop it entirely.
$ws.Rows.Item(17).Delete()

# Make this the active sheet/selection, matching the new tabSelected/selection.
$ws.Activate()
$ws.Range("C1:H4").Select()
